# "2024 Day 22 Done"
# Fill in Day 21 and Day 22 results on the "2024" sheet, update the running
# totals for every day after those two days were completed, record the
# per-day status codes for days 21/22 on the "Overall" sheet, and leave the
# selection/active-sheet state the way the author left it (on "Overall",
# cell A13 selected) after having last worked on "2024" cell C24.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("2024")
$ws2 = $wb.Worksheets.Item("Overall")

# "2024" has protected cells (B/C/E/F are unlocked inputs, the rest are
# locked formula cells) - lift protection while we write, then restore it.
$ws1.Unprotect()

# Updated "Done Both" (B) / "Done 1 Only" (C) counts for days 1-20 (rows 2-21)
$dayTotals = @{
  2  = @(232441, 17025)
  3  = @(165422, 40678)
  4  = @(150210, 16299)
  5  = @(121145, 10409)
  6  = @(99755,  12940)
  7  = @(74369,  25670)
  8  = @(78125,  4129)
  9  = @(67521,  2831)
  10 = @(57210,  11054)
  11 = @(58624,  1269)
  12 = @(53738,  8634)
  13 = @(40028,  10826)
  14 = @(41744,  5029)
  15 = @(39701,  4512)
  16 = @(29882,  9245)
  17 = @(25315,  5611)
  18 = @(22036,  9549)
  19 = @(27190,  661)
  20 = @(24440,  2639)
  21 = @(18707,  3710)
}

foreach ($row in $dayTotals.Keys) {
  $vals = $dayTotals[$row]
  $ws1.Range("B$row").Value = $vals[0]
  $ws1.Range("C$row").Value = $vals[1]
}

# Day 21 (row 22) and Day 22 (row 23) newly solved - fill in raw results
$ws1.Range("B22").Value = 9353
$ws1.Range("C22").Value = 3829
$ws1.Range("E22").Value = 6458

$ws1.Range("B23").Value = 9657
$ws1.Range("C23").Value = 3570
$ws1.Range("E23").Value = 10291
$ws1.Range("F23").Value = 8994

# Re-assert the (unchanged) formulas on rows 22/23 so the engine refreshes
# their cached results now that the row has real inputs (picking up the
# newly-typed B/C/E/F and the refreshed $B$2/$D$2 totals above).
$ws1.Range("D22").Formula = '=IF(ISBLANK(B22),"",B22+C22)'
$ws1.Range("G22").Formula = '=IF(D22="","",E22/D22)'
$ws1.Range("H22").Formula = '=IF(ISBLANK(C22),"",F22/B22)'
$ws1.Range("I22").Formula = '=IF(ISBLANK(E22),"",E22/$D$2)'
$ws1.Range("J22").Formula = '=IF(ISBLANK(F22),"",F22/$B$2)'

$ws1.Range("D23").Formula = '=IF(ISBLANK(B23),"",B23+C23)'
$ws1.Range("G23").Formula = '=IF(D23="","",E23/D23)'
$ws1.Range("H23").Formula = '=IF(ISBLANK(C23),"",F23/B23)'
$ws1.Range("I23").Formula = '=IF(ISBLANK(E23),"",E23/$D$2)'
$ws1.Range("J23").Formula = '=IF(ISBLANK(F23),"",F23/$B$2)'

$ws1.Protect()

# Record day 21 / day 22 status codes on the Overall sheet (row 13 = 2024)
$ws2.Range("CD13").Value = "s"
$ws2.Range("CE13").Value = "p"
$ws2.Range("CF13").Value = "m"
$ws2.Range("CG13").Value = "m"
$ws2.Range("CH13").Value = "s"
$ws2.Range("CI13").Value = "s"
$ws2.Range("CJ13").Value = "s"
$ws2.Range("CK13").Value = "s"

# Leave the selection on "2024" at C24 (where the author last clicked before
# switching sheets), then make "Overall" the active sheet with A13 selected.
$ws1.Activate()
$ws1.Range("C24").Select()

$ws2.Activate()
$ws2.Range("A13").Select()
